# Trade #14 closed at 2026-02-16 22:58:42 - base_strategy UP +0.000%
# Append the new trade row (row 15) to both the "All Trades" and
# "base_strategy" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 15

    $ws.Range("A$row").Value = 14
    # Leading apostrophe forces text so the date-like string isn't
    # auto-converted into a date serial number (matches the other rows).
    $ws.Range("B$row").Value = "'2026-02-16"
    $ws.Range("C$row").Value = "22:58:42"
    $ws.Range("D$row").Value = "base_strategy"
    $ws.Range("E$row").Value = "UP"
    $ws.Range("F$row").Value = 0.5
    # Lone apostrophe -> empty text cell (matches empty Exit Price cell).
    $ws.Range("G$row").Value = "'"
    $ws.Range("H$row").Value = "OPEN"
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 100
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = "'"
    $ws.Range("Q$row").Value = 0
}
